# Apply the "fix yolk issue" changes:
#  - Rename the "Yolk_Ave"/"Yolk_SD" columns (N, O) to "Embryo_Ave"/"Embryo_SD"
#  - Update the related data-validation rules (header-name guards + value-range
#    guards) to reference the new names
#  - Widen the Conductivity upper bound from 718 to 781
#  - Re-range the Month validation to a warning-style model-extrapolation
#    check (4-8) instead of a hard whole-number stop (1-12)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells: rename Yolk_* to Embryo_* -----------------------------
$ws.Range("N1").Value = "Embryo_Ave"
$ws.Range("O1").Value = "Embryo_SD"

# --- Conductivity (F2:F1048576): widen upper bound 718 -> 781 ------------
$dvF = $ws.Range("F2:F1048576").Validation
$dvF.Formula1 = 274
$dvF.Formula2 = 781

# --- Embryo_Ave values (N2:N1048576): update warning text ----------------
$dvN = $ws.Range("N2:N1048576").Validation
$dvN.Formula1 = 0.434
$dvN.Formula2 = 4.371
$dvN.ErrorTitle = "Warning"
$dvN.ErrorMessage = "Embryo_Ave entered falls outside of the range of the data using to train the model (0.434 to 4.371 mm). The model will be forced to extrapolate when making a prediction. Would you like to proceed?"
$dvN.InputTitle = "Average Yolk Diameter"
$dvN.InputMessage = "mm"
$dvN.AlertStyle = 2

# --- Embryo_SD values (O2:O1048576): update warning text -----------------
$dvO = $ws.Range("O2:O1048576").Validation
$dvO.Formula1 = 0.005
$dvO.Formula2 = 1.377
$dvO.ErrorTitle = "Warning"
$dvO.ErrorMessage = "Embryo_SD entered falls outside of the range of the data using to train the model (0.005 to 1.377 mm). The model will be forced to extrapolate when making a prediction. Would you like to proceed?"
$dvO.InputTitle = "SD of Yolk Diameter"
$dvO.InputMessage = "mm"
$dvO.AlertStyle = 2

# --- Month (C2:C1048576): switch to warning-style extrapolation check ----
$dvC = $ws.Range("C2:C1048576").Validation
$dvC.Formula1 = 4
$dvC.Formula2 = 8
$dvC.ErrorTitle = "Entry Error"
$dvC.ErrorMessage = "Month entered falls outside of the range of the data using to train the models (4,5,6,7,8). The model will be forced to extrapolate when making a prediction. Would you like to proceed?"
$dvC.InputTitle = "Month"
$dvC.InputMessage = "1 to 12"
$dvC.AlertStyle = 2

# --- Header-name guard on N1 (Average Yolk Diameter -> Average Embryo Diameter)
$dvN1 = $ws.Range("N1").Validation
$dvN1.Formula1 = """Embryo_Ave"""
$dvN1.ErrorTitle = "Error"
$dvN1.ErrorMessage = "Required variable names cannot be adjusted"
$dvN1.InputTitle = "Average Embryo Diameter"
$dvN1.InputMessage = "mm"

# --- Header-name guard on O1 (SD of Yolk Diameter -> SD of Embryo Diameter)
$dvO1 = $ws.Range("O1").Validation
$dvO1.Formula1 = """Embryo_SD"""
$dvO1.ErrorTitle = "Error"
$dvO1.ErrorMessage = "Required variable names cannot be adjusted"
$dvO1.InputTitle = "SD of Embryo Diameter"
$dvO1.InputMessage = "mm"
